# Insert two new rows at 228-229, pushing existing rows 228..290 down to 230..292.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("228:229").Insert()

# Row 228 (new): Sutil De Gase / Primera, Peru, $/caja 24 kilos
$ws.Cells.Item(228, 1).Value2 = 1
$ws.Cells.Item(228, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(228, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(228, 4).Value2 = 44841
$ws.Cells.Item(228, 5).Value2 = 15
$ws.Cells.Item(228, 6).Value2 = "Fruta"
$ws.Cells.Item(228, 7).Value2 = 100102
$ws.Cells.Item(228, 8).Value2 = "Cítricos"
$ws.Cells.Item(228, 9).Value2 = 100102003
$ws.Cells.Item(228, 10).Value2 = "Limón"
$ws.Cells.Item(228, 11).Value2 = "Sutil De Gase"
$ws.Cells.Item(228, 12).Value2 = "Primera"
$ws.Cells.Item(228, 13).Value2 = 250
$ws.Cells.Item(228, 14).Value2 = 32000
$ws.Cells.Item(228, 15).Value2 = 33000
$ws.Cells.Item(228, 16).Value2 = 32500
$ws.Cells.Item(228, 17).Value2 = "`$/caja 24 kilos"
$ws.Cells.Item(228, 18).Value2 = "Perú"
$ws.Cells.Item(228, 19).Value2 = 1354
$ws.Cells.Item(228, 20).Value2 = 24

# Row 229 (new): Tahití / Primera, Peru, $/caja 24 kilos
$ws.Cells.Item(229, 1).Value2 = 1
$ws.Cells.Item(229, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(229, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(229, 4).Value2 = 44841
$ws.Cells.Item(229, 5).Value2 = 15
$ws.Cells.Item(229, 6).Value2 = "Fruta"
$ws.Cells.Item(229, 7).Value2 = 100102
$ws.Cells.Item(229, 8).Value2 = "Cítricos"
$ws.Cells.Item(229, 9).Value2 = 100102003
$ws.Cells.Item(229, 10).Value2 = "Limón"
$ws.Cells.Item(229, 11).Value2 = "Tahití"
$ws.Cells.Item(229, 12).Value2 = "Primera"
$ws.Cells.Item(229, 13).Value2 = 300
$ws.Cells.Item(229, 14).Value2 = 29000
$ws.Cells.Item(229, 15).Value2 = 30000
$ws.Cells.Item(229, 16).Value2 = 29500
$ws.Cells.Item(229, 17).Value2 = "`$/caja 24 kilos"
$ws.Cells.Item(229, 18).Value2 = "Perú"
$ws.Cells.Item(229, 19).Value2 = 1229
$ws.Cells.Item(229, 20).Value2 = 24
